$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 397.2857
$ws.Range("J80").Value = 492.33334
$ws.Range("L80").Value = 1477.00002
$ws.Range("N80").Value = -3473.00002
$ws.Range("H83").Value = 397.2857
$ws.Range("J83").Value = 492.33334
$ws.Range("L83").Value = 4431.00006
$ws.Range("N83").Value = -14415.00006
$ws.Range("H88").Value = 3783323
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 3783323
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 3783323
$ws.Range("M88").ClearContents()
$ws.Range("N88").Value = -3784135
$ws.Range("H91").Value = 3783323
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 3783323
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 3783323
$ws.Range("M91").ClearContents()
$ws.Range("N91").Value = -3786131
$ws.Range("H132").Value = 2767.6667
$ws.Range("I132").Value = 2724.75
$ws.Range("K132").Value = 8174.25
$ws.Range("M132").Value = -5644.25
$ws.Range("H137").Value = 38463900
$ws.Range("J137").Value = 844.5
$ws.Range("L137").Value = 2533.5
$ws.Range("N137").Value = -7633.5
$ws.Range("H138").Value = 2074.532
$ws.Range("I138").Value = 1559.4706
$ws.Range("K138").Value = 4678.4118
$ws.Range("M138").Value = 461.5882000000001
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3923.8223
$ws.Range("I32").Value = 2237.6765
$ws.Range("K32").Value = 2237.6765
$ws.Range("M32").Value = -1950.6765
$ws.Range("H46").Value = 5168.4
$ws.Range("I46").Value = 4000
$ws.Range("J46").Value = 5460.5
$ws.Range("K46").Value = 4000
$ws.Range("L46").Value = 5460.5
$ws.Range("M46").Value = -3681
$ws.Range("N46").Value = -6098.5
$ws.Range("H61").Value = 3563.56
$ws.Range("I61").Value = 2024.55
$ws.Range("K61").Value = 2024.55
$ws.Range("M61").Value = -1812.55
$ws.Range("H136").Value = 3563.56
$ws.Range("I136").Value = 2024.55
$ws.Range("K136").Value = 6073.65
$ws.Range("M136").Value = -3523.65
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 997.5833
$ws.Range("J80").Value = 791.6667
$ws.Range("L80").Value = 791.6667
$ws.Range("N80").Value = -2787.6667
$ws.Range("H83").Value = 997.5833
$ws.Range("J83").Value = 791.6667
$ws.Range("L83").Value = 3958.3335
$ws.Range("N83").Value = -13942.3335
$ws.Range("H105").Value = 1187.091
$ws.Range("I105").Value = 926.3333
$ws.Range("K105").Value = 926.3333
$ws.Range("M105").Value = 820.6667
$ws.Range("H134").Value = 588027.9399999999
$ws.Range("I134").Value = 811823.4
$ws.Range("K134").Value = 2435470.2
$ws.Range("M134").Value = -2432935.2
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2242.3215
$ws.Range("I31").Value = 1853.7826
$ws.Range("K31").Value = 1853.7826
$ws.Range("M31").Value = -1558.7826
$ws.Range("H34").Value = 2242.3215
$ws.Range("I34").Value = 1853.7826
$ws.Range("K34").Value = 1853.7826
$ws.Range("M34").Value = -1651.7826
$ws.Range("H62").Value = 11663
$ws.Range("I62").Value = 5000
$ws.Range("J62").Value = 14994.5
$ws.Range("K62").Value = 5000
$ws.Range("L62").Value = 14994.5
$ws.Range("M62").Value = -4376
$ws.Range("N62").Value = -16242.5
$ws.Range("H65").Value = 11663
$ws.Range("I65").Value = 5000
$ws.Range("J65").Value = 14994.5
$ws.Range("K65").Value = 25000
$ws.Range("L65").Value = 74972.5
$ws.Range("M65").Value = -21880
$ws.Range("N65").Value = -81212.5
$ws.Range("H69").Value = 15000
$ws.Range("I69").Value = 15000
$ws.Range("K69").Value = 15000
$ws.Range("M69").Value = -14251
$ws.Range("H72").Value = 15000
$ws.Range("I72").Value = 15000
$ws.Range("K72").Value = 45000
$ws.Range("M72").Value = -41256
$ws.Range("H105").Value = 7837.393
$ws.Range("I105").Value = 9993.286
$ws.Range("K105").Value = 9993.286
$ws.Range("M105").Value = -8246.286
$ws.Range("H132").Value = 2516.25
$ws.Range("I132").Value = 2306.6052
$ws.Range("J132").Value = 6499.5
$ws.Range("K132").Value = 6919.8156
$ws.Range("L132").Value = 19498.5
$ws.Range("M132").Value = -4389.8156
$ws.Range("N132").Value = -24558.5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 244.6
$ws.Range("I11").Value = 268.25
$ws.Range("K11").Value = 804.75
$ws.Range("M11").Value = -664.75
$ws.Range("H12").Value = 399.42105
$ws.Range("I12").Value = 61
$ws.Range("K12").Value = 183
$ws.Range("M12").Value = -10
$ws.Range("H86").Value = 940
$ws.Range("I86").Value = 666.6667
$ws.Range("J86").Value = 1350
$ws.Range("K86").Value = 2000.0001
$ws.Range("L86").Value = 4050
$ws.Range("M86").Value = -814.0001
$ws.Range("N86").Value = -6422
$ws.Range("H89").Value = 940
$ws.Range("I89").Value = 666.6667
$ws.Range("J89").Value = 1350
$ws.Range("K89").Value = 6000.0003
$ws.Range("L89").Value = 12150
$ws.Range("M89").Value = -72.0002999999997
$ws.Range("N89").Value = -24006
$ws.Range("H125").Value = 2995.25
$ws.Range("I125").Value = 2993.6667
$ws.Range("K125").Value = 8981.000100000001
$ws.Range("M125").Value = -4061.000100000001
$ws.Range("H129").Value = 6671225
$ws.Range("I129").Value = 5336.857
$ws.Range("K129").Value = 16010.571
$ws.Range("M129").Value = -11010.571
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3718.2173
$ws.Range("I102").Value = 3801.4
$ws.Range("K102").Value = 3801.4
$ws.Range("M102").Value = -2179.4
$ws.Range("H122").Value = 2940.682
$ws.Range("J122").Value = 2479.7144
$ws.Range("L122").Value = 7439.1432
$ws.Range("N122").Value = -12339.1432
$ws.Range("H132").Value = 2713.4167
$ws.Range("I132").Value = 2696.4546
$ws.Range("K132").Value = 8089.3638
$ws.Range("M132").Value = -5559.3638
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 52635936
$ws.Range("I7").Value = 111113224
$ws.Range("K7").Value = 111113224
$ws.Range("M7").Value = -111113112
$ws.Range("H126").Value = 52635936
$ws.Range("I126").Value = 111113224
$ws.Range("K126").Value = 333339672
$ws.Range("M126").Value = -333337202
$ws.Range("H132").Value = 3514
$ws.Range("I132").Value = 3514
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 10542
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -8012
$ws.Range("N132").ClearContents()
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 105179
$ws.Range("J62").Value = 16409.8
$ws.Range("L62").Value = 16409.8
$ws.Range("N62").Value = -17657.8
$ws.Range("H65").Value = 105179
$ws.Range("J65").Value = 16409.8
$ws.Range("L65").Value = 82049
$ws.Range("N65").Value = -88289
$ws.Range("H70").Value = 33000
$ws.Range("I70").Value = 16000
$ws.Range("J70").Value = 50000
$ws.Range("K70").Value = 16000
$ws.Range("L70").Value = 50000
$ws.Range("M70").Value = -15685
$ws.Range("N70").Value = -50630
$ws.Range("H73").Value = 33000
$ws.Range("I73").Value = 16000
$ws.Range("J73").Value = 50000
$ws.Range("K73").Value = 16000
$ws.Range("L73").Value = 50000
$ws.Range("M73").Value = -14908
$ws.Range("N73").Value = -52184
$ws.Range("H122").Value = 112504130
$ws.Range("I122").Value = 83338830
$ws.Range("K122").Value = 250016490
$ws.Range("M122").Value = -250014040
